$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Sheet "VENTAS POR GRUPO": M19 (PORCELANATO for TAMAYO VILLACIS EDWIN XAVIER) 0 -> 33.7
$wsGrupo.Range("M19").Value = 33.7

# Sheet "VENTAS POR GRUPO": M24 summary text "7 de 22" -> "8 de 22"
$wsGrupo.Range("M24").Value = "8 de 22"

# Sheet "VENTA MENSUAL": F19 (julio for TAMAYO VILLACIS EDWIN XAVIER) 0 -> 33.7
$wsMensual.Range("F19").Value = 33.7

# Sheet "VENTA MENSUAL": F24 total julio column 33382.14 -> 33415.84
$wsMensual.Range("F24").Value = 33415.84

# Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO row (16) VENTA, POR CUMPLIR and CUMPLIMIENTO updated
$wsCumpl.Range("D16").Value = 29878.86
$wsCumpl.Range("E16").Value = 8877.68
$wsCumpl.Range("F16").Value = 0.7709372405276632

# Sheet "CUMPLIMIENTO MENSUAL": TOTAL row (19) VENTA, POR CUMPLIR and CUMPLIMIENTO updated
$wsCumpl.Range("D19").Value = 33415.84
$wsCumpl.Range("E19").Value = 24807.16386304604
$wsCumpl.Range("F19").Value = 0.5739284781424501

# Column E width auto-adjusted slightly (23 -> 22) as a side effect of the data update
$wsCumpl.Columns.Item(5).ColumnWidth = 21.14
